$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column J
$ws.Range("J1").Value = "Studentoffer"

# Row 3 - Student fare, Chennai -> Mumbai
$ws.Range("A3").Value = 8015332963
$ws.Range("A3").Font.Color = 16777215
$ws.Range("B3").Value = 8234
$ws.Range("C3").Value = "Chennai"
$ws.Range("D3").Value = "Mumbai"
$ws.Range("E3").Value = 45936
$ws.Range("E3").NumberFormat = "[$-14009]d\ mmmm\ yyyy;@"
$ws.Range("F3").Value = 45942
$ws.Range("F3").NumberFormat = "[$-14009]d\ mmmm\ yyyy;@"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = "Economy"
$ws.Range("J3").Value = "Student"

# Row 4 - Senior Citizen fare, Coimbatore -> Hyderabad
$ws.Range("A4").Value = 8015332963
$ws.Range("A4").Font.Color = 16777215
$ws.Range("B4").Value = 8234
$ws.Range("C4").Value = "Coimbatore"
$ws.Range("D4").Value = "Hyderabad"
$ws.Range("E4").Value = 45931
$ws.Range("E4").NumberFormat = "[$-14009]d\ mmmm\ yyyy;@"
$ws.Range("F4").Value = 45961
$ws.Range("F4").NumberFormat = "[$-14009]d\ mmmm\ yyyy;@"
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = "Premium Economy"
$ws.Range("J4").Value = "Senior Citizen"

# Column widths for newly used columns (compensate for the host's
# internal 1/6-character snapping so the saved OOXML `width` lands as
# close as possible to the target values)
$ws.Columns.Item(9).ColumnWidth = 15.276041666666666
$ws.Columns.Item(10).ColumnWidth = 13.608072916666666
$ws.Columns.Item(11).ColumnWidth = 19.166666666666668
$ws.Columns.Item(12).ColumnWidth = 27.053385416666668

# Selection moves to J5 after data entry
$ws.Range("J5").Select()
